# Apply the c19stats-ac-kennzahlen.xlsx update:
#  - Fill in the newly computed 7-day-moving-average columns
#    (H, J, M, O, Q) for the existing rows 159-161.
#  - Append three new days of data as rows 162-164
#    (dates 2020-08-08, 2020-08-09, 2020-08-10), whose H/J/M/O/Q
#    moving-average columns have not been computed yet (left blank,
#    matching the source rows they were copied from before calc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backfill the moving-average columns for rows 159-161 ---

$ws.Range("H159").Value = 4.142857142857143
$ws.Range("J159").Value = 7.58876576486639
$ws.Range("M159").Value = 8.857142857142858
$ws.Range("O159").Value = 0
$ws.Range("Q159").Value = 4.714285714285714

$ws.Range("H160").Value = 2.428571428571428
$ws.Range("J160").Value = 5.402847679270441
$ws.Range("M160").Value = 8.857142857142858
$ws.Range("O160").Value = 0
$ws.Range("Q160").Value = 6.428571428571429

$ws.Range("H161").Value = 2.428571428571428
$ws.Range("J161").Value = 5.402847679270441
$ws.Range("M161").Value = 8.857142857142858
$ws.Range("O161").Value = 0
$ws.Range("Q161").Value = 6.428571428571429

# --- Append the three new data rows (162-164) ---
# Copy the date cell's style (s="2", the YYYY-MM-DD HH:MM:SS date format)
# from the last existing row down to the new rows, then overwrite the
# values on every cell that should hold data.

$ws.Range("A161").Copy($ws.Range("A162"))
$ws.Range("A161").Copy($ws.Range("A163"))
$ws.Range("A161").Copy($ws.Range("A164"))

# Row 162 : 2020-08-08
$ws.Range("A162").Value = 44051
$ws.Range("B162").Value = 2157
$ws.Range("C162").Value = 1060
$ws.Range("D162").Value = 100
$ws.Range("E162").Value = 1975
$ws.Range("F162").Value = 82
$ws.Range("G162").Value = 0
$ws.Range("I162").Value = 0
$ws.Range("K162").Value = 0
$ws.Range("L162").Value = 0
$ws.Range("N162").Value = 0
$ws.Range("P162").Value = 0
$ws.Range("R162").Value = 388.3232966973616
$ws.Range("S162").Value = 428.4905812919395
$ws.Range("T162").Value = 11.72285552591155
$ws.Range("U162").Value = 11.16181937655838

# Row 163 : 2020-08-09
$ws.Range("A163").Value = 44052
$ws.Range("B163").Value = 2170
$ws.Range("C163").Value = 1068
$ws.Range("D163").Value = 100
$ws.Range("E163").Value = 1999
$ws.Range("F163").Value = 71
$ws.Range("G163").Value = -11
$ws.Range("I163").Value = -13.41463414634146
$ws.Range("K163").Value = 13
$ws.Range("L163").Value = 8
$ws.Range("N163").Value = 0
$ws.Range("P163").Value = 24
$ws.Range("R163").Value = 390.6636781795433
$ws.Range("S163").Value = 431.7244724715013
$ws.Range("T163").Value = 11.31861912846633
$ws.Range("U163").Value = 11.16181937655838

# Row 164 : 2020-08-10
$ws.Range("A164").Value = 44053
$ws.Range("B164").Value = 2170
$ws.Range("C164").Value = 1068
$ws.Range("D164").Value = 100
$ws.Range("E164").Value = 1999
$ws.Range("F164").Value = 71
$ws.Range("G164").Value = 0
$ws.Range("I164").Value = 0
$ws.Range("K164").Value = 0
$ws.Range("L164").Value = 0
$ws.Range("N164").Value = 0
$ws.Range("P164").Value = 0
$ws.Range("R164").Value = 390.6636781795433
$ws.Range("S164").Value = 431.7244724715013
$ws.Range("T164").Value = 11.31861912846633
$ws.Range("U164").Value = 11.16181937655838
